$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.017.52"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "3.389.35"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'572.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Value = "'141.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("E9").Value = "  +2.54%  "
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").Value = "'0.387"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.45%  "
$ws.Range("D12").Value = "3.968.18"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "'27.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.393.15"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000171"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").Value = "61.130.38"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").Value = "'6.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.38%  "
$ws.Range("D19").Value = "'13.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("D20").Value = "'8.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("D21").Value = "'384.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").Value = "'75.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.91%  "
$ws.Range("D23").Value = "'0.552"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").Value = "3.524.81"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  +3.14%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'7.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").Value = "'2.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").Value = "'1.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.73%  "
$ws.Range("D34").Value = "'23.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.03%  "
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("D36").Value = "'166.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "3.422.96"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("E39").Value = "  -2.30%  "
$ws.Range("E40").Value = "  -1.06%  "
$ws.Range("D41").Value = "'26.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.94%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "'4.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").Value = "2.454.82"
$ws.Range("E47").Value = "  -3.21%  "
$ws.Range("D48").Value = "'22.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("E49").Value = "  -1.93%  "
$ws.Range("E50").Value = "  +11.12%  "
$ws.Range("E51").Value = "  -0.85%  "
